$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row updates ---
$ws.Range("C1").Value = "CNPJ do Vendedor"

# --- Row 2: wrap existing values into Python-list-like string representations ---
$ws.Range("B2").Value = "['374.17']"
$ws.Range("C2").Value = "['36882195000279']"
$ws.Range("D2").Value = "['BROTA COMPANY COMERCIO DE PLANTAS LTDA']"
$ws.Range("E2").Value = "['10000000000']"
$ws.Range("F2").Value = "['Lira da Hashtag']"
$ws.Range("G2").Value = "[[('b.box led', '389.00'), ('Tomilho Serpilho', '10.00'), ('Oregano Bravo-Europeu', '10.00'), ('Manjericao Italiano', '10.00'), ('Coentro Portugues', '10.00'), ('Salsa Hortense', '10.00'), ('Alface Baby-Leaf', '10.00')]]"

# --- Row 3: new row of data ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "['849.00']"
$ws.Range("C3").Value = "['60409075055054']"
$ws.Range("D3").Value = "['NESTLE BRASIL LTDA']"
$ws.Range("E3").Value = "['11122233344']"
$ws.Range("F3").Value = "['Lira da Hashtag']"
$ws.Range("G3").Value = "[[('ESSENZA Mini C30 Metal 110V', '534.82'), ('NOMAD Travel Mug Festve Medium', '104.35'), ('Pack YEP 2021 50 Caps OL', '130.00')]]"

# --- Row 4: new row of data ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "['4500']"
$ws.Range("C4").Value = "['30000000000101']"
$ws.Range("D4").Value = "['LIRA BOLADO NO XML']"
$ws.Range("E4").Value = "['26344392000108']"
$ws.Range("F4").Value = "['HASHTAG TREINAMENTOS LTDA']"
$ws.Range("G4").Value = "['Nota referente aos serviços realizados no mês de Novembro de 2021.']"
